$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.724.13'
$ws.Range('E2').Value = '  +2.01%  '
$ws.Range('D3').Value = '2.116.11'
$ws.Range('E3').Value = '  +10.37%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '333.99'
$ws.Range('E5').Value = '  +4.31%  '
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5226'
$ws.Range('E7').Value = '  +3.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4419'
$ws.Range('E8').Value = '  +8.50%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09119'
$ws.Range('E9').Value = '  +9.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '46.43'
$ws.Range('E10').Value = '  +9.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.185'
$ws.Range('E11').Value = '  +6.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '25.22'
$ws.Range('E12').Value = '  +4.91%  '
$ws.Range('D13').Value = '2.114.10'
$ws.Range('E13').Value = '  +10.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.802'
$ws.Range('E14').Value = '  +5.92%  '
$ws.Range('E15').Value = '  +7.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '98.43'
$ws.Range('E16').Value = '  +6.36%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001142'
$ws.Range('E17').Value = '  +4.22%  '
$ws.Range('B18').Value = 'BinanceUSD'
$ws.Range('C18').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.002'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.28'
$ws.Range('E20').Value = '  +4.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.412'
$ws.Range('E21').Value = '  +7.91%  '
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').Value = '30.832.98'
$ws.Range('E23').Value = '  +2.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.07'
$ws.Range('E24').Value = '  +6.24%  '
$ws.Range('D25').Value = '2.361.91'
$ws.Range('E25').Value = '  +10.99%  '
$ws.Range('E26').Value = '  +2.95%  '
$ws.Range('E27').Value = '  +5.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.574'
$ws.Range('E28').Value = '  +13.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '163.66'
$ws.Range('E29').Value = '  +0.51%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.82'
$ws.Range('E30').Value = '  +3.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.190'
$ws.Range('E31').Value = '  +4.17%  '
$ws.Range('E32').Value = '  +2.37%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.252'
$ws.Range('E33').Value = '  +5.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.912'
$ws.Range('E34').Value = '  +3.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.526'
$ws.Range('E35').Value = '  +27.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02594'
$ws.Range('E36').Value = '  +5.58%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.596'
$ws.Range('E37').Value = '  +4.43%  '
$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '9.646'
$ws.Range('E38').Value = '  +12.06%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06786'
$ws.Range('E39').Value = '  +5.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.81'
$ws.Range('E40').Value = '  +12.32%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2269'
$ws.Range('E41').Value = '  +5.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6810'
$ws.Range('E42').Value = '  +4.39%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.255'
$ws.Range('E43').Value = '  +3.63%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.32'
$ws.Range('E44').Value = '  +7.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9994'
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6365'
$ws.Range('E46').Value = '  +4.76%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.258'
$ws.Range('E47').Value = '  +3.08%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.676'
$ws.Range('E48').Value = '  +1.44%  '
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.292'
$ws.Range('E49').Value = '  +6.80%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '83.33'
$ws.Range('E50').Value = '  +5.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '119.13'
$ws.Range('E51').Value = '  -2.61%  '
